$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (nameId|String) - shifts star/accumulatedAtk/prefabAddress/orderIndex/meetWeight right
$ws.Columns("B").Insert()

# Insert a new column F (spriteName|String) - after prefabAddress (now column E), shifts orderIndex/meetWeight right
$ws.Columns("F").Insert()

# Header row
$ws.Range("B1").Value = "nameId|String"
$ws.Range("F1").Value = "spriteName|String"

# Data rows: B = PetName_000X, F = PetPortrait_0001 (write row-by-row to match shared-string insertion order)
$ws.Range("B2").Value = "PetName_0001"
$ws.Range("F2").Value = "PetPortrait_0001"

$ws.Range("B3").Value = "PetName_0002"
$ws.Range("F3").Value = "PetPortrait_0001"

$ws.Range("B4").Value = "PetName_0003"
$ws.Range("F4").Value = "PetPortrait_0001"

$ws.Range("B5").Value = "PetName_0004"
$ws.Range("F5").Value = "PetPortrait_0001"

$ws.Range("B6").Value = "PetName_0005"
$ws.Range("F6").Value = "PetPortrait_0001"

$ws.Range("B7").Value = "PetName_0006"
$ws.Range("F7").Value = "PetPortrait_0001"

# Reset the selection to the default top-left cell (matches the author's final view state)
$ws.Range("A1").Select() | Out-Null
